$p = $ppt.ActivePresentation

# The deck currently ends with slide 8 ("Questions", SlideID 260).
# We insert two brand-new "Title and Content" slides right before it,
# which pushes "Questions" down to the end (position 10) while the two
# new slides take positions 8 and 9 - matching the sldIdLst reordering
# (264, 265, 260) in the target revision.

$graphSlide = $p.Slides.Add(8, 2)
$graphSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Identifying Leaks with Graphs"
$graphSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Construct a directed graph`rNodes: DOM elements and JavaScript objects`rEdges: References DOMElement---->JSObj and JSObj---->DOMElement`rAdd nodes and edges when DOM element refers a JS object and vice versa`rUse Graph algorithms to identify all the cycles in the graph. If cycle exists then there is a circular reference that may result in memory leak`r"

$contSlide = $p.Slides.Add(9, 2)
$contSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Identifying Leaks with Graphs cont."
$contSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Will identify most of the leaks resulting from Circular References, and Closures with Circular References`rMost of the leaks in the existing applications are of this kind`r`r"
